# Applies the weekly odds-sheet update described in the commit:
#  - refresh several odds figures on row 3 (Bolivia match)
#  - refresh two odds figures on row 4 (Brazil match)
#  - remove the Paraguay match row (old row 7), shifting the three
#    following matches (USA/Peru/USA rows) up by one row

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 3 odds updates -------------------------------------------------
$ws.Range("G3").Value  = 1.57
$ws.Range("I3").Value  = 5.25
$ws.Range("L3").Value  = 5
$ws.Range("W3").Value  = 9
$ws.Range("Z3").Value  = 12
$ws.Range("AD3").Value = 8.5
$ws.Range("AE3").Value = 15
$ws.Range("AK3").Value = 41
$ws.Range("AM3").Value = 151
$ws.Range("AQ3").Value = 21
$ws.Range("AX3").Value = 26
$ws.Range("AY3").Value = 29

# --- Row 4 odds updates -------------------------------------------------
$ws.Range("Q4").Value = 2.07
$ws.Range("R4").Value = 1.83

# --- Remove the Paraguay - Primera Division row (old row 7) -------------
# This shifts the USA/Peru/USA rows that followed it up by one position
# (old rows 8,9,10 become new rows 7,8,9) and updates the sheet
# dimension from A1:BD10 to A1:BD9 automatically.
$ws.Rows.Item(7).Delete()
